$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "description" column header to "descriptions" and populate each
# row with its own JSON array of bullet-point descriptions (replacing the
# single shared Lorem-ipsum placeholder text). Cells must be written in the
# same order the strings were originally authored so the shared-string table
# is rebuilt in the same sequence as the target workbook.

$ws.Range("B3").Value = "[`n""Design, implement, secure, and maintain AWS environments for multiple applications on the platform"",`n""Implement Jenkins jobs and pipelines to analyze code and deploy code artifacts onto active servers"",`n""Perform cost optimization on AWS environments to meet monetary requirements"",`n""Provide development guidance on best practices for numerous programming languages to numerous colleagues"",`n""Write Ansible playbooks and CloudFormation templates to automate AWS infrastructure""`n]"

$ws.Range("B1").Value = "descriptions"

$ws.Range("B5").Value = "[`n""Coordinate API communication across numerous teams working on a unified product for the business"",`n""Design and implement backend logic for <a href='https://www.verizon.com/5g/home' target='_blank'>Verizon's 5G Home Customer Qualification</a>""`n]"

$ws.Range("B4").Value = "[`n""Design and implement backend logic including API specifications and database elements in Java and Python"",`n""Translate business requirements into technical requirements so they can be written as code"",`n""Write application code to generate tiles embedded with data for serving on <a href='https://www.verizon.com/coverage-map/' target='_blank'>Verizon's 5G Coverage Map</a>""`n]"

$ws.Range("B6").Value = "[`n""Act as gate keeper to multiple backend Java repositories"",`n""Perform maintenance and optimizations on existing Java code written by colleagues"",`n""Provide Java code guidance to the team when required"",`n""Write Java and unit tests to complete technical requirements""`n]"

$ws.Range("B7").Value = "[`n""Convert legacy code to match current Java standards"",`n""Perform data entry and correct data bugs in multiple databases""`n]"

$ws.Range("B8").Value = "[`n""Created an expense reimbursement sample application"",`n""Created an inventory management sample application"",`n""Created an internal blogging application for technical trainers""`n]"

$ws.Range("B2").Value = "[`n""Implement full stack code with pair programming style against  vigorous standards"",`n""Implement multiple types of tests to maintain 100% average"",`n""Interact with GCP products to perform analysis, corrections, monitoring, logging, and validation""`n]"

# Update the saved cursor/selection position recorded in the sheet view.
[void]$ws.Range("E13").Select()
